# Update insert_normunit for commercial building type compatibility.
# - Adds a new "NumBldgs" worksheet (after "NumStor") holding a lookup
#   table of BldgType -> numbldgs counts.
# - The rows are entered in their "natural" order first, then four rows
#   (Fin, Gro, Lib, Rel) are inserted afterwards to keep the sheet sorted
#   alphabetically -- this reproduces the shared-string insertion order
#   seen in the target workbook (new unique strings are appended to the
#   shared string table in first-seen order, not sheet order).
# - NumStor's selection/view is moved off the top of the sheet and the
#   new NumBldgs sheet becomes the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new worksheet right after NumStor ---------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "NumBldgs"

# --- header row ---------------------------------------------------------
$ws2.Range("A1").Value = "BldgType"
$ws2.Range("B1").Value = "numbldgs"

# --- data rows, in the order they were originally typed (pre-sort, and
#     before the Fin/Gro/Lib/Rel/SUn/WRf rows existed) -------------------
$initialData = @(
    @("Dmo", 2),
    @("MFm", 24),
    @("SFm", 2),
    @("Asm", 1),
    @("ECC", 1),
    @("EPr", 1),
    @("ERC", 1),
    @("ESe", 1),
    @("EUn", 1),
    @("Hsp", 1),
    @("Htl", 1),
    @("MBT", 1),
    @("MLI", 1),
    @("Mtl", 1),
    @("Nrs", 1),
    @("OfL", 1),
    @("OfS", 1),
    @("RFF", 1),
    @("RSD", 1),
    @("Rt3", 1),
    @("RtL", 1),
    @("RtS", 1),
    @("SCn", 1)
)

$row = 2
foreach ($item in $initialData) {
    $ws2.Cells.Item($row, 1).Value = $item[0]
    $ws2.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# --- insert the four rows that keep the sheet alphabetically sorted ----
# (Fin after EUn, Gro after Fin, Lib after Htl, Rel after OfS)
$ws2.Rows.Item(11).Insert()
$ws2.Cells.Item(11, 1).Value = "Fin"
$ws2.Cells.Item(11, 2).Value = 1

$ws2.Rows.Item(12).Insert()
$ws2.Cells.Item(12, 1).Value = "Gro"
$ws2.Cells.Item(12, 2).Value = 1

$ws2.Rows.Item(15).Insert()
$ws2.Cells.Item(15, 1).Value = "Lib"
$ws2.Cells.Item(15, 2).Value = 1

$ws2.Rows.Item(22).Insert()
$ws2.Cells.Item(22, 1).Value = "Rel"
$ws2.Cells.Item(22, 2).Value = 1

# --- finally append SUn and WRf as brand-new rows at the bottom --------
$ws2.Cells.Item(29, 1).Value = "SUn"
$ws2.Cells.Item(29, 2).Value = 1
$ws2.Cells.Item(30, 1).Value = "WRf"
$ws2.Cells.Item(30, 2).Value = 1

# --- view/selection bookkeeping -----------------------------------------
# NumStor: scroll/select away from the top (A100 region / N114).
$ws1.Range("N114").Select()

# NumBldgs: select B5:B30 (mirrors the author's on-screen selection) and
# make this sheet the active tab, like in the saved workbook.
$ws2.Range("B5:B30").Select()
$ws2.Activate()
